$d = $word.ActiveDocument

# Locate the three bullet paragraphs we need to touch by their text,
# rather than relying on fixed indices.
$pTrends = $null
$pRepeat = $null
$pMaps = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "Create trends for the different seasons*") { $pTrends = $p }
    elseif ($t -like "Repeat the comparison between the DBs*") { $pRepeat = $p }
    elseif ($t -like "Create composite maps for all RSTs combined*") { $pMaps = $p }
}

# 1) Strike through the "Create trends..." bullet (now done/cancelled).
$pTrends.Range.Font.StrikeThrough = 1

# 2) Strike through the "Repeat the comparison..." bullet (now done/cancelled).
$pRepeat.Range.Font.StrikeThrough = 1

# 3) Split "Create composite maps..." bullet into two runs, moving the
#    _GoBack bookmark to sit between "Create composite " and
#    "maps for all RSTs combined, for 00z and 12z."
$prefix = "Create composite "
$splitPos = $pMaps.Range.Start + [int]$prefix.Length

$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitRange)

Write-Output "done"
